$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tr = $sh.TextFrame.TextRange

# Paragraph 2 (lvl 1): "We came very close..." -> "We would like to have working navigation..."
$tr.Paragraphs(2,1).Runs(1,1).Text = "We would like to have working navigation in a good looking  navbar."

# Paragraph 3 (lvl 1): "We would also like to have working navigation..." -> "Photo references..."
$tr.Paragraphs(3,1).Runs(1,1).Text = "Photo references are returned in the JSON but photos were not Implemented in the front end."

# Paragraph 4 (lvl 1): "Photo references..." -> split into two runs:
#   "(pull up postman to look at JSON responses & " + "photo responses)"
$run4 = $tr.Paragraphs(4,1).Runs(1,1)
$run4.Text = "(pull up postman to look at JSON responses & "
$run4.InsertAfter("photo responses)") | Out-Null
